$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B -- this shifts the existing
#    B:H data (and their column widths) right to C:I.
$ws.Columns("B:B").Insert()

# 2. Copy the formatting of the (now shifted) cells onto the new column B
#    cells so the style indexes match exactly, then set the values
#    ("# of Diffs" header + per-row diff counts).

# Header B2 <- style copied from C2 (old header style, bold/filled).
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B2").Value = "# of Diffs"

# Styled (grey-highlight) rows: copy from D3, which is the shifted cell
# that already carries that highlighted style.
$ws.Range("D3").Copy()
$styledRows = @(3,4,5,6,7,8,10)
foreach ($r in $styledRows) {
    $ws.Range("B$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$diffCounts = @{3=3.0; 4=3.0; 5=3.0; 6=3.0; 7=3.0; 8=3.0; 9=0.0; 10=1.0}
foreach ($r in 3..10) {
    $ws.Range("B$r").Value = $diffCounts[$r]
}

# The ID column (now C) never carries the grey-highlight style in the
# target layout, even though a couple of rows happened to inherit it
# from the old B column via the Insert() shift -- strip it back off.
$ws.Range("C3:C10").ClearFormats()

# 3. Column B width ("# of Diffs" column).
$ws.Columns("B:B").ColumnWidth = 13.166666666666666

# 4. Fix up the autoFilter range (Insert left the old, stale B2:H10 ref
#    in place) and the _FilterDatabase defined name so both cover the new
#    column (B2:I10 / B2:I2 instead of B2:H10 / B2:H2).
$ws.AutoFilterMode = $false
$ws.Range("B2:I10").AutoFilter()

$wb.Names("_xlnm._FilterDatabase").RefersTo = "=full!`$B`$2:`$I`$2"
